$d = $word.ActiveDocument
$r = $d.Content
$r.Find.Execute("36×68=", $true, $false, $false, $false, $false, $true, 1, $false, "88×35=", 2) | Out-Null
$r.Find.Execute("30×32=", $true, $false, $false, $false, $false, $true, 1, $false, "27×18=", 2) | Out-Null
$r.Find.Execute("80×63=", $true, $false, $false, $false, $false, $true, 1, $false, "22×90=", 2) | Out-Null
$r.Find.Execute("19×73=", $true, $false, $false, $false, $false, $true, 1, $false, "14×52=", 2) | Out-Null
$r.Find.Execute("41×13=", $true, $false, $false, $false, $false, $true, 1, $false, "19×62=", 2) | Out-Null
$r.Find.Execute("15×41=", $true, $false, $false, $false, $false, $true, 1, $false, "75×77=", 2) | Out-Null
$r.Find.Execute("40×27=", $true, $false, $false, $false, $false, $true, 1, $false, "48×43=", 2) | Out-Null
$r.Find.Execute("85×50=", $true, $false, $false, $false, $false, $true, 1, $false, "47×63=", 2) | Out-Null
$r.Find.Execute("34×53=", $true, $false, $false, $false, $false, $true, 1, $false, "15×32=", 2) | Out-Null
$r.Find.Execute("50×43=", $true, $false, $false, $false, $false, $true, 1, $false, "92×13=", 2) | Out-Null
$r.Find.Execute("77×75=", $true, $false, $false, $false, $false, $true, 1, $false, "38×36=", 2) | Out-Null
$r.Find.Execute("58×24=", $true, $false, $false, $false, $false, $true, 1, $false, "82×17=", 2) | Out-Null
$r.Find.Execute("22×10=", $true, $false, $false, $false, $false, $true, 1, $false, "58×64=", 2) | Out-Null
$r.Find.Execute("49×57=", $true, $false, $false, $false, $false, $true, 1, $false, "62×23=", 2) | Out-Null
$r.Find.Execute("10×18=", $true, $false, $false, $false, $false, $true, 1, $false, "89×24=", 2) | Out-Null
$r.Find.Execute("29×43=", $true, $false, $false, $false, $false, $true, 1, $false, "37×91=", 2) | Out-Null
$r.Find.Execute("73×11=", $true, $false, $false, $false, $false, $true, 1, $false, "98×27=", 2) | Out-Null
$r.Find.Execute("16×52=", $true, $false, $false, $false, $false, $true, 1, $false, "39×24=", 2) | Out-Null
$r.Find.Execute("67×24=", $true, $false, $false, $false, $false, $true, 1, $false, "12×74=", 2) | Out-Null
$r.Find.Execute("80×79=", $true, $false, $false, $false, $false, $true, 1, $false, "90×39=", 2) | Out-Null
$r.Find.Execute("23×32=", $true, $false, $false, $false, $false, $true, 1, $false, "73×83=", 2) | Out-Null
$r.Find.Execute("76×45=", $true, $false, $false, $false, $false, $true, 1, $false, "15×47=", 2) | Out-Null
$r.Find.Execute("55×69=", $true, $false, $false, $false, $false, $true, 1, $false, "97×77=", 2) | Out-Null
$r.Find.Execute("97×27=", $true, $false, $false, $false, $false, $true, 1, $false, "88×60=", 2) | Out-Null
$r.Find.Execute("64×82=", $true, $false, $false, $false, $false, $true, 1, $false, "83×80=", 2) | Out-Null
$r.Find.Execute("59×84=", $true, $false, $false, $false, $false, $true, 1, $false, "83×19=", 2) | Out-Null
$r.Find.Execute("72×50=", $true, $false, $false, $false, $false, $true, 1, $false, "46×92=", 2) | Out-Null
$r.Find.Execute("81×75=", $true, $false, $false, $false, $false, $true, 1, $false, "18×31=", 2) | Out-Null
$r.Find.Execute("16×91=", $true, $false, $false, $false, $false, $true, 1, $false, "81×45=", 2) | Out-Null
$r.Find.Execute("67×32=", $true, $false, $false, $false, $false, $true, 1, $false, "52×86=", 2) | Out-Null
$r.Find.Execute("10×56=", $true, $false, $false, $false, $false, $true, 1, $false, "60×74=", 2) | Out-Null
$r.Find.Execute("66×36=", $true, $false, $false, $false, $false, $true, 1, $false, "31×37=", 2) | Out-Null
$r.Find.Execute("40×100=", $true, $false, $false, $false, $false, $true, 1, $false, "28×55=", 2) | Out-Null
$r.Find.Execute("91×40=", $true, $false, $false, $false, $false, $true, 1, $false, "77×50=", 2) | Out-Null
$r.Find.Execute("29×49=", $true, $false, $false, $false, $false, $true, 1, $false, "52×38=", 2) | Out-Null
$r.Find.Execute("63×88=", $true, $false, $false, $false, $false, $true, 1, $false, "26×52=", 2) | Out-Null
$r.Find.Execute("86×92=", $true, $false, $false, $false, $false, $true, 1, $false, "95×33=", 2) | Out-Null
$r.Find.Execute("66×69=", $true, $false, $false, $false, $false, $true, 1, $false, "34×91=", 2) | Out-Null
$r.Find.Execute("77×97=", $true, $false, $false, $false, $false, $true, 1, $false, "46×48=", 2) | Out-Null
$r.Find.Execute("69×15=", $true, $false, $false, $false, $false, $true, 1, $false, "54×59=", 2) | Out-Null
$r.Find.Execute("29×87=", $true, $false, $false, $false, $false, $true, 1, $false, "55×52=", 2) | Out-Null
$r.Find.Execute("29×19=", $true, $false, $false, $false, $false, $true, 1, $false, "97×48=", 2) | Out-Null
$r.Find.Execute("55×11=", $true, $false, $false, $false, $false, $true, 1, $false, "70×65=", 2) | Out-Null
$r.Find.Execute("15×66=", $true, $false, $false, $false, $false, $true, 1, $false, "97×29=", 2) | Out-Null
$r.Find.Execute("73×12=", $true, $false, $false, $false, $false, $true, 1, $false, "37×12=", 2) | Out-Null
$r.Find.Execute("25×19=", $true, $false, $false, $false, $false, $true, 1, $false, "60×57=", 2) | Out-Null
$r.Find.Execute("28×76=", $true, $false, $false, $false, $false, $true, 1, $false, "17×45=", 2) | Out-Null
$r.Find.Execute("49×56=", $true, $false, $false, $false, $false, $true, 1, $false, "32×88=", 2) | Out-Null
$r.Find.Execute("72×20=", $true, $false, $false, $false, $false, $true, 1, $false, "72×61=", 2) | Out-Null
$r.Find.Execute("11×76=", $true, $false, $false, $false, $false, $true, 1, $false, "99×27=", 2) | Out-Null
$r.Find.Execute("61×91=", $true, $false, $false, $false, $false, $true, 1, $false, "65×83=", 2) | Out-Null
$r.Find.Execute("75×76=", $true, $false, $false, $false, $false, $true, 1, $false, "86×10=", 2) | Out-Null
$r.Find.Execute("81×20=", $true, $false, $false, $false, $false, $true, 1, $false, "23×62=", 2) | Out-Null
$r.Find.Execute("91×54=", $true, $false, $false, $false, $false, $true, 1, $false, "50×42=", 2) | Out-Null
$r.Find.Execute("60×34=", $true, $false, $false, $false, $false, $true, 1, $false, "86×12=", 2) | Out-Null
$r.Find.Execute("25×36=", $true, $false, $false, $false, $false, $true, 1, $false, "42×67=", 2) | Out-Null
$r.Find.Execute("68×46=", $true, $false, $false, $false, $false, $true, 1, $false, "63×21=", 2) | Out-Null
$r.Find.Execute("46×96=", $true, $false, $false, $false, $false, $true, 1, $false, "41×38=", 2) | Out-Null
$r.Find.Execute("44×61=", $true, $false, $false, $false, $false, $true, 1, $false, "76×35=", 2) | Out-Null
$r.Find.Execute("94×37=", $true, $false, $false, $false, $false, $true, 1, $false, "33×33=", 2) | Out-Null
$r.Find.Execute("38×28=", $true, $false, $false, $false, $false, $true, 1, $false, "49×19=", 2) | Out-Null
$r.Find.Execute("82×72=", $true, $false, $false, $false, $false, $true, 1, $false, "89×88=", 2) | Out-Null
$r.Find.Execute("34×50=", $true, $false, $false, $false, $false, $true, 1, $false, "45×84=", 2) | Out-Null
$r.Find.Execute("60×37=", $true, $false, $false, $false, $false, $true, 1, $false, "34×27=", 2) | Out-Null
$r.Find.Execute("89×43=", $true, $false, $false, $false, $false, $true, 1, $false, "71×65=", 2) | Out-Null
$r.Find.Execute("70×79=", $true, $false, $false, $false, $false, $true, 1, $false, "54×81=", 2) | Out-Null
$r.Find.Execute("48×68=", $true, $false, $false, $false, $false, $true, 1, $false, "23×29=", 2) | Out-Null
$r.Find.Execute("69×91=", $true, $false, $false, $false, $false, $true, 1, $false, "58×90=", 2) | Out-Null
$r.Find.Execute("28×66=", $true, $false, $false, $false, $false, $true, 1, $false, "97×43=", 2) | Out-Null
$r.Find.Execute("16×87=", $true, $false, $false, $false, $false, $true, 1, $false, "26×18=", 2) | Out-Null
$r.Find.Execute("68×80=", $true, $false, $false, $false, $false, $true, 1, $false, "25×32=", 2) | Out-Null
$r.Find.Execute("89×16=", $true, $false, $false, $false, $false, $true, 1, $false, "31×37=", 2) | Out-Null
$r.Find.Execute("76×93=", $true, $false, $false, $false, $false, $true, 1, $false, "43×39=", 2) | Out-Null
$r.Find.Execute("22×67=", $true, $false, $false, $false, $false, $true, 1, $false, "21×55=", 2) | Out-Null
$r.Find.Execute("89×99=", $true, $false, $false, $false, $false, $true, 1, $false, "45×54=", 2) | Out-Null
$r.Find.Execute("21×35=", $true, $false, $false, $false, $false, $true, 1, $false, "49×83=", 2) | Out-Null
$r.Find.Execute("43×45=", $true, $false, $false, $false, $false, $true, 1, $false, "78×41=", 2) | Out-Null
$r.Find.Execute("22×78=", $true, $false, $false, $false, $false, $true, 1, $false, "31×81=", 2) | Out-Null
$r.Find.Execute("39×22=", $true, $false, $false, $false, $false, $true, 1, $false, "14×60=", 2) | Out-Null
$r.Find.Execute("22×83=", $true, $false, $false, $false, $false, $true, 1, $false, "64×73=", 2) | Out-Null
$r.Find.Execute("33×48=", $true, $false, $false, $false, $false, $true, 1, $false, "33×71=", 2) | Out-Null
$r.Find.Execute("84×52=", $true, $false, $false, $false, $false, $true, 1, $false, "63×22=", 2) | Out-Null
$r.Find.Execute("64×34=", $true, $false, $false, $false, $false, $true, 1, $false, "90×37=", 2) | Out-Null
$r.Find.Execute("53×94=", $true, $false, $false, $false, $false, $true, 1, $false, "89×98=", 2) | Out-Null
$r.Find.Execute("72×39=", $true, $false, $false, $false, $false, $true, 1, $false, "27×39=", 2) | Out-Null
$r.Find.Execute("37×15=", $true, $false, $false, $false, $false, $true, 1, $false, "59×12=", 2) | Out-Null
$r.Find.Execute("81×57=", $true, $false, $false, $false, $false, $true, 1, $false, "98×73=", 2) | Out-Null
$r.Find.Execute("92×34=", $true, $false, $false, $false, $false, $true, 1, $false, "31×73=", 2) | Out-Null
$r.Find.Execute("29×66=", $true, $false, $false, $false, $false, $true, 1, $false, "90×85=", 2) | Out-Null
$r.Find.Execute("98×34=", $true, $false, $false, $false, $false, $true, 1, $false, "54×24=", 2) | Out-Null
$r.Find.Execute("39×13=", $true, $false, $false, $false, $false, $true, 1, $false, "99×18=", 2) | Out-Null
$r.Find.Execute("20×67=", $true, $false, $false, $false, $false, $true, 1, $false, "86×98=", 2) | Out-Null
$r.Find.Execute("45×93=", $true, $false, $false, $false, $false, $true, 1, $false, "38×78=", 2) | Out-Null
$r.Find.Execute("38×74=", $true, $false, $false, $false, $false, $true, 1, $false, "18×46=", 2) | Out-Null
$r.Find.Execute("24×98=", $true, $false, $false, $false, $false, $true, 1, $false, "44×92=", 2) | Out-Null
$r.Find.Execute("69×56=", $true, $false, $false, $false, $false, $true, 1, $false, "91×63=", 2) | Out-Null
$r.Find.Execute("87×10=", $true, $false, $false, $false, $false, $true, 1, $false, "87×30=", 2) | Out-Null
$r.Find.Execute("91×87=", $true, $false, $false, $false, $false, $true, 1, $false, "43×54=", 2) | Out-Null
$r.Find.Execute("46×39=", $true, $false, $false, $false, $false, $true, 1, $false, "13×89=", 2) | Out-Null
$r.Find.Execute("71×73=", $true, $false, $false, $false, $false, $true, 1, $false, "96×65=", 2) | Out-Null
